{"js": "// Update the date line and the 25 \"a\u00f7b=c, d\" answer cells in the single\n// 20x5 table (populated rows are table rows 0, 4, 8, 12, 16 (0-based); the\n// rest are blank spacer rows).\n//\n// Each cell/paragraph is targeted by position (paragraph index / row+col),\n// not by searching for its old text: several \"after\" values coincide with\n// other entries' \"before\" values (e.g. row 0's new text equals row 1's old\n// text), so a plain sequential search-and-replace pass could re-match text\n// it had just written. Updating in place via insertText(..., replace) on\n// the existing paragraph/range also keeps the original run formatting\n// (font, size, alignment) intact.\n\nconst body = context.document.body;\n\n// --- Date heading (first paragraph in the body) ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2025-01-22 Wednesday\", Word.InsertLocation.replace);\n\n// --- Table of division problems ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst rowUpdates = [\n  [0, [\"88\u00f75=17, 3\", \"49\u00f79=5, 4\", \"56\u00f78=7, 0\", \"74\u00f77=10, 4\", \"51\u00f77=7, 2\"]],\n  [4, [\"41\u00f73=13, 2\", \"39\u00f73=13, 0\", \"91\u00f74=22, 3\", \"54\u00f76=9, 0\", \"66\u00f79=7, 3\"]],\n  [8, [\"79\u00f73=26, 1\", \"50\u00f72=25, 0\", \"25\u00f73=8, 1\", \"82\u00f74=20, 2\", \"13\u00f76=2, 1\"]],\n  [12, [\"30\u00f79=3, 3\", \"72\u00f78=9, 0\", \"78\u00f76=13, 0\", \"29\u00f78=3, 5\", \"70\u00f78=8, 6\"]],\n  [16, [\"88\u00f75=17, 3\", \"73\u00f75=14, 3\", \"16\u00f78=2, 0\", \"37\u00f73=12, 1\", \"83\u00f72=41, 1\"]],\n];\n\nfor (const [rowIndex, values] of rowUpdates) {\n  for (let col = 0; col < values.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    const cellParagraphs = cell.body.paragraphs;\n    cellParagraphs.load(\"items\");\n    await context.sync();\n    cellParagraphs.items[0].insertText(values[col], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 \"a\u00f7b=c, d\" answer cells in the single\n# 20x5 table (5 populated rows at table rows 1,5,9,13,17; the rest are\n# blank spacer rows).\n$d = $word.ActiveDocument\n\n# Date paragraph (first paragraph in the document body).\n$d.Paragraphs.Item(1).Range.Text = \"2025-01-22 Wednesday\"\n\n$t = $d.Tables.Item(1)\n\n$rowsData = @(\n    @{ Row = 1;  Values = @(\"88\u00f75=17, 3\", \"49\u00f79=5, 4\", \"56\u00f78=7, 0\", \"74\u00f77=10, 4\", \"51\u00f77=7, 2\") },\n    @{ Row = 5;  Values = @(\"41\u00f73=13, 2\", \"39\u00f73=13, 0\", \"91\u00f74=22, 3\", \"54\u00f76=9, 0\", \"66\u00f79=7, 3\") },\n    @{ Row = 9;  Values = @(\"79\u00f73=26, 1\", \"50\u00f72=25, 0\", \"25\u00f73=8, 1\", \"82\u00f74=20, 2\", \"13\u00f76=2, 1\") },\n    @{ Row = 13; Values = @(\"30\u00f79=3, 3\", \"72\u00f78=9, 0\", \"78\u00f76=13, 0\", \"29\u00f78=3, 5\", \"70\u00f78=8, 6\") },\n    @{ Row = 17; Values = @(\"88\u00f75=17, 3\", \"73\u00f75=14, 3\", \"16\u00f78=2, 0\", \"37\u00f73=12, 1\", \"83\u00f72=41, 1\") }\n)\n\nforeach ($rowInfo in $rowsData) {\n    $r = $rowInfo.Row\n    $vals = $rowInfo.Values\n    for ($c = 1; $c -le 5; $c++) {\n        $t.Cell($r, $c).Range.Text = $vals[$c - 1]\n    }\n}\n"}
